# Updating scripts for Pagination
#
# - Widen column A on the TestData sheet slightly.
# - Clear the old ad-hoc test-case note that used to live in A10.
# - Populate the previously-blank rows 11-15 with the new pagination
#   configuration variables / values used by the automation scripts.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestData")

# Widen column A (27.86 -> 31.14 characters) to fit the new, longer
# variable names added below (e.g. "paginationDashboardNamePageSize").
$ws.Columns.Item(1).ColumnWidth = 31.14

# Row 10 no longer carries the old manual test-case description.
$ws.Range("A10").Value2 = $null

# New pagination configuration rows.
$ws.Range("A11").Value2 = "paginationDashboardName"
$ws.Range("B11").Value2 = "Automation_Pivot_Pagination_Dashboard"

$ws.Range("A12").Value2 = "paginationInsightName"
$ws.Range("B12").Value2 = "Automation_Insight"

$ws.Range("A13").Value2 = "paginationDashboardNamePageSize"
$ws.Range("B13").Value2 = "Automation_Pagination_Dashboard_PageSize"

$ws.Range("A14").Value2 = "PageSize"
$ws.Range("B14").Value2 = 10

$ws.Range("A15").Value2 = "DefaultPageSize"
$ws.Range("B15").Value2 = 10000
$ws.Range("B15").NumberFormat = "#,##0"
